$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update codigo values
$ws.Range("A2").Value = 1234456
$ws.Range("E2").Value = 12.32
$ws.Range("A3").Value = 2324567

# Update image path text for the perfume row (shared string content change)
$ws.Range("C3").Value = "..\..\Imagenes\2324567-perfume.jpg"

# Delete the last row (cerveza entry) entirely
$ws.Rows("4:4").Delete()
